$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 130.57143
$ws.Range("I9").Value = 166.33333
$ws.Range("K9").Value = 166.33333
$ws.Range("M9").Value = 2.666670000000011
$ws.Range("H17").Value = 753.8873
$ws.Range("J17").Value = 786.9552
$ws.Range("L17").Value = 2360.8656
$ws.Range("N17").Value = -2696.8656
$ws.Range("H38").Value = 2412.75
$ws.Range("I38").Value = 2400.2856
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 7200.8568
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = -6828.8568
$ws.Range("N38").Value = -8244
$ws.Range("H74").Value = 4726.4443
$ws.Range("I74").Value = 3257
$ws.Range("K74").Value = 3257
$ws.Range("M74").Value = -2321
$ws.Range("H77").Value = 4726.4443
$ws.Range("I77").Value = 3257
$ws.Range("K77").Value = 16285
$ws.Range("M77").Value = -11605
$ws.Range("H100").Value = 4277
$ws.Range("I100").Value = 3627
$ws.Range("K100").Value = 3627
$ws.Range("M100").Value = -3086
$ws.Range("H141").Value = 2089.1365
$ws.Range("I141").Value = 2224.8
$ws.Range("K141").Value = 6674.400000000001
$ws.Range("M141").Value = -1494.400000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50722
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52496
$ws.Range("H102").Value = 2089.889
$ws.Range("I102").Value = 1670
$ws.Range("J102").Value = 2299.8333
$ws.Range("K102").Value = 1670
$ws.Range("L102").Value = 2299.8333
$ws.Range("M102").Value = -48
$ws.Range("N102").Value = -5543.8333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 11732
$ws.Range("I10").Value = 3990
$ws.Range("J10").Value = 15603
$ws.Range("K10").Value = 3990
$ws.Range("L10").Value = 15603
$ws.Range("M10").Value = -3850
$ws.Range("N10").Value = -15883
$ws.Range("H20").Value = 1732.6333
$ws.Range("I20").Value = 1213.238
$ws.Range("K20").Value = 1213.238
$ws.Range("M20").Value = -966.2380000000001
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H99").Value = 4156.077
$ws.Range("I99").Value = 4200.4443
$ws.Range("K99").Value = 4200.4443
$ws.Range("M99").Value = -2702.4443
$ws.Range("H107").Value = 3117.1428
$ws.Range("I107").Value = 3095
$ws.Range("K107").Value = 3095
$ws.Range("M107").Value = -1175
$ws.Range("H134").Value = 25651162
$ws.Range("I134").Value = 10980.333
$ws.Range("K134").Value = 32940.999
$ws.Range("M134").Value = -30405.999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4268
$ws.Range("J22").Value = 4575.2
$ws.Range("L22").Value = 4575.2
$ws.Range("N22").Value = -5275.2
$ws.Range("H80").Value = 20127.75
$ws.Range("J80").Value = 20127.75
$ws.Range("L80").Value = 20127.75
$ws.Range("N80").Value = -22373.75
$ws.Range("H82").Value = 80090.5
$ws.Range("J82").Value = 80090.5
$ws.Range("L82").Value = 80090.5
$ws.Range("N82").Value = -80812.5
$ws.Range("H83").Value = 20127.75
$ws.Range("J83").Value = 20127.75
$ws.Range("L83").Value = 60383.25
$ws.Range("N83").Value = -71615.25
$ws.Range("H85").Value = 80090.5
$ws.Range("J85").Value = 80090.5
$ws.Range("L85").Value = 80090.5
$ws.Range("N85").Value = -82586.5
$ws.Range("H87").Value = 8999
$ws.Range("J87").Value = 8999
$ws.Range("L87").Value = 8999
$ws.Range("N87").Value = -11371
$ws.Range("H88").Value = 40114
$ws.Range("J88").Value = 40114
$ws.Range("L88").Value = 40114
$ws.Range("N88").Value = -40926
$ws.Range("H90").Value = 8999
$ws.Range("J90").Value = 8999
$ws.Range("L90").Value = 26997
$ws.Range("N90").Value = -38853
$ws.Range("H91").Value = 40114
$ws.Range("J91").Value = 40114
$ws.Range("L91").Value = 40114
$ws.Range("N91").Value = -42922
$ws.Range("H122").Value = 1379.25
$ws.Range("I122").Value = 1012
$ws.Range("K122").Value = 3036
$ws.Range("M122").Value = -586
$ws.Range("H134").Value = 7146979
$ws.Range("I134").Value = 3803.7144
$ws.Range("J134").Value = 14290154
$ws.Range("K134").Value = 11411.1432
$ws.Range("L134").Value = 42870462
$ws.Range("M134").Value = -8876.143199999999
$ws.Range("N134").Value = -42875532

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1008.381
$ws.Range("I2").Value = 1613.75
$ws.Range("J2").Value = 635.8461
$ws.Range("K2").Value = 9682.5
$ws.Range("L2").Value = 3815.0766
$ws.Range("M2").Value = -9569.5
$ws.Range("N2").Value = -4041.0766
$ws.Range("H29").Value = 929.25
$ws.Range("J29").Value = 2117
$ws.Range("L29").Value = 6351
$ws.Range("N29").Value = -6905
$ws.Range("H40").Value = 156.8077
$ws.Range("J40").Value = 232.18182
$ws.Range("L40").Value = 928.72728
$ws.Range("N40").Value = -1066.72728
$ws.Range("H68").Value = 1899.7273
$ws.Range("I68").Value = 1120.75
$ws.Range("J68").Value = 3977
$ws.Range("K68").Value = 3362.25
$ws.Range("L68").Value = 11931
$ws.Range("M68").Value = -2551.25
$ws.Range("N68").Value = -13553
$ws.Range("H71").Value = 1899.7273
$ws.Range("I71").Value = 1120.75
$ws.Range("J71").Value = 3977
$ws.Range("K71").Value = 10086.75
$ws.Range("L71").Value = 35793
$ws.Range("M71").Value = -6030.75
$ws.Range("N71").Value = -43905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H122").Value = 1620.4736
$ws.Range("I122").Value = 1822.4722
$ws.Range("J122").Value = 1274.1904
$ws.Range("K122").Value = 5467.4166
$ws.Range("L122").Value = 3822.5712
$ws.Range("M122").Value = -3017.4166
$ws.Range("N122").Value = -8722.5712
$ws.Range("H132").Value = 2617.8064
$ws.Range("I132").Value = 2521.9524
$ws.Range("J132").Value = 2819.1
$ws.Range("K132").Value = 7565.8572
$ws.Range("L132").Value = 8457.299999999999
$ws.Range("M132").Value = -5035.8572
$ws.Range("N132").Value = -13517.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 732.25
$ws.Range("I16").Value = 732.25
$ws.Range("K16").Value = 732.25
$ws.Range("M16").Value = -562.25
$ws.Range("H22").Value = 1011.3333
$ws.Range("I22").Value = 423
$ws.Range("J22").Value = 1599.6666
$ws.Range("K22").Value = 423
$ws.Range("L22").Value = 1599.6666
$ws.Range("M22").Value = -128
$ws.Range("N22").Value = -2189.6666
$ws.Range("H27").Value = 1011.3333
$ws.Range("I27").Value = 423
$ws.Range("J27").Value = 1599.6666
$ws.Range("K27").Value = 423
$ws.Range("L27").Value = 1599.6666
$ws.Range("M27").Value = -316
$ws.Range("N27").Value = -1813.6666
$ws.Range("H64").Value = 62500
$ws.Range("J64").Value = 62500
$ws.Range("L64").Value = 62500
$ws.Range("N64").Value = -62950
$ws.Range("H67").Value = 62500
$ws.Range("J67").Value = 62500
$ws.Range("L67").Value = 62500
$ws.Range("N67").Value = -64060
$ws.Range("H69").Value = 55081.5
$ws.Range("J69").Value = 55081.5
$ws.Range("L69").Value = 55081.5
$ws.Range("N69").Value = -56703.5
$ws.Range("H72").Value = 55081.5
$ws.Range("J72").Value = 55081.5
$ws.Range("L72").Value = 165244.5
$ws.Range("N72").Value = -173356.5
$ws.Range("H114").Value = 43499.668
$ws.Range("I114").Value = 39999
$ws.Range("J114").Value = 45250
$ws.Range("K114").Value = 39999
$ws.Range("L114").Value = 45250
$ws.Range("M114").Value = -35660
$ws.Range("N114").Value = -53928

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H122").Value = 2689.6843
$ws.Range("I122").Value = 2800.6875
$ws.Range("J122").Value = 2097.6667
$ws.Range("K122").Value = 8402.0625
$ws.Range("L122").Value = 6293.000100000001
$ws.Range("M122").Value = -5952.0625
$ws.Range("N122").Value = -11193.0001
$ws.Range("H132").Value = 1537.8889
$ws.Range("I132").Value = 1260.56
$ws.Range("K132").Value = 3781.68
$ws.Range("M132").Value = -1251.68
